$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 182, shifting existing rows 182:199 down to 184:201.
$ws.Rows("182:183").Insert()

# Row 182 is a new weekly entry for "Primera" quality, same data as the (now shifted)
# row 184 except for the date, which moves forward to 44449.
$ws.Range("A182").Value = 1
$ws.Range("B182").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C182").Value = "Arica y Parinacota"
$ws.Range("D182").Value = 44449
$ws.Range("E182").Value = 15
$ws.Range("F182").Value = 100112032
$ws.Range("G182").Value = "Zapallo italiano"
$ws.Range("H182").Value = "Huracán"
$ws.Range("I182").Value = "Primera"
$ws.Range("J182").Value = 120
$ws.Range("K182").Value = 11000
$ws.Range("L182").Value = 12000
$ws.Range("M182").Value = 11500
$ws.Range("N182").Value = '$/caja 70 unidades'
$ws.Range("O182").Value = "Región de Arica y Parinacota"
$ws.Range("P182").Value = 164
$ws.Range("Q182").Value = 70
$ws.Range("R182").Value = "Hortaliza"

# Row 183 is a new weekly entry for "Segunda" quality, same data as the (now shifted)
# row 185 except for the date, which moves forward to 44449.
$ws.Range("A183").Value = 1
$ws.Range("B183").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C183").Value = "Arica y Parinacota"
$ws.Range("D183").Value = 44449
$ws.Range("E183").Value = 15
$ws.Range("F183").Value = 100112032
$ws.Range("G183").Value = "Zapallo italiano"
$ws.Range("H183").Value = "Huracán"
$ws.Range("I183").Value = "Segunda"
$ws.Range("J183").Value = 120
$ws.Range("K183").Value = 9000
$ws.Range("L183").Value = 10000
$ws.Range("M183").Value = 9500
$ws.Range("N183").Value = '$/caja 100 unidades'
$ws.Range("O183").Value = "Región de Arica y Parinacota"
$ws.Range("P183").Value = 95
$ws.Range("Q183").Value = 100
$ws.Range("R183").Value = "Hortaliza"
